$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with changed odds (rows 4, 6, 7) ---
# Row 4
$ws.Range("G4").Value = 1.81
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 4.75
$ws.Range("J4").Value = 2.63
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 5.5
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
$ws.Range("U4").Value = 2.25
$ws.Range("V4").Value = 1.57
$ws.Range("X4").Value = 7.5
$ws.Range("Y4").Value = 9.5
$ws.Range("Z4").Value = 15
$ws.Range("AA4").Value = 19
$ws.Range("AC4").Value = 6.5
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 21
$ws.Range("AH4").Value = 9.5
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 17
$ws.Range("AK4").Value = 51
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 11
$ws.Range("AP4").Value = 26
$ws.Range("AR4").Value = 67
$ws.Range("AX4").Value = 6
$ws.Range("AY4").Value = 29
$ws.Range("BA4").Value = 101
$ws.Range("BB4").Value = 151

# Row 6
$ws.Range("P6").Value = 4.02
$ws.Range("T6").Value = 3.21

# Row 7
$ws.Range("O7").Value = 1.19
$ws.Range("P7").Value = 4.33
$ws.Range("U7").Value = 1.87
$ws.Range("V7").Value = 1.87

# --- Insert two new match rows before the old row 8 (Ukraine match shifts to row 10) ---
$ws.Range("8:9").Insert()

# Row 8 - new match data
$ws.Range("A8").Value = "U5ehv3Vb"
$ws.Range("B8").Value = "25/10/2024"
$ws.Range("C8").Value = "11:35"
$ws.Range("D8").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E8").Value = "Al Qadisiya"
$ws.Range("F8").Value = "Damac"
$ws.Range("G8").Value = 1.65
$ws.Range("H8").Value = 3.9
$ws.Range("I8").Value = 4.75
$ws.Range("J8").Value = 2.2
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 4.75
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 9
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 3.75
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 1.36
$ws.Range("T8").Value = 3
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.83
$ws.Range("W8").Value = 7.5
$ws.Range("X8").Value = 8
$ws.Range("Y8").Value = 8.5
$ws.Range("Z8").Value = 12
$ws.Range("AA8").Value = 13
$ws.Range("AB8").Value = 26
$ws.Range("AC8").Value = 11
$ws.Range("AD8").Value = 7.5
$ws.Range("AE8").Value = 17
$ws.Range("AF8").Value = 51
$ws.Range("AG8").Value = 600
$ws.Range("AH8").Value = 13
$ws.Range("AI8").Value = 26
$ws.Range("AJ8").Value = 17
$ws.Range("AK8").Value = 51
$ws.Range("AL8").Value = 41
$ws.Range("AM8").Value = 41
$ws.Range("AN8").Value = 3.75
$ws.Range("AO8").Value = 8.5
$ws.Range("AP8").Value = 19
$ws.Range("AQ8").Value = 26
$ws.Range("AR8").Value = 51
$ws.Range("AS8").Value = 126
$ws.Range("AT8").Value = 3
$ws.Range("AU8").Value = 8.5
$ws.Range("AV8").Value = 51
$ws.Range("AW8").Value = 81
$ws.Range("AX8").Value = 6.5
$ws.Range("AY8").Value = 26
$ws.Range("AZ8").Value = 34
$ws.Range("BA8").Value = 81
$ws.Range("BB8").Value = 101
$ws.Range("BC8").Value = 400
$ws.Range("BD8").Value = 81

# Row 9 - new match data
$ws.Range("A9").Value = "lnKaZbQH"
$ws.Range("B9").Value = "25/10/2024"
$ws.Range("C9").Value = "11:30"
$ws.Range("D9").Value = "SERBIA - SUPER LIGA"
$ws.Range("E9").Value = "IMT Novi Beograd"
$ws.Range("F9").Value = "Tekstilac Odzaci"
$ws.Range("G9").Value = 1.5
$ws.Range("H9").Value = 4.05
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 2.05
$ws.Range("K9").Value = 2.2
$ws.Range("L9").Value = 5.9
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 7.3
$ws.Range("O9").Value = 1.31
$ws.Range("P9").Value = 3.2
$ws.Range("Q9").Value = 1.93
$ws.Range("R9").Value = 1.82
$ws.Range("S9").Value = 1.42
$ws.Range("T9").Value = 2.67
$ws.Range("U9").Value = 2.05
$ws.Range("V9").Value = 1.7
$ws.Range("W9").Value = 6
$ws.Range("X9").Value = 6.4
$ws.Range("Y9").Value = 8.25
$ws.Range("Z9").Value = 10
$ws.Range("AA9").Value = 13
$ws.Range("AB9").Value = 32
$ws.Range("AC9").Value = 7.3
$ws.Range("AD9").Value = 7.8
$ws.Range("AE9").Value = 20
$ws.Range("AF9").Value = 110
$ws.Range("AG9").Value = 1000
$ws.Range("AH9").Value = 14
$ws.Range("AI9").Value = 35
$ws.Range("AJ9").Value = 19.5
$ws.Range("AK9").Value = 120
$ws.Range("AL9").Value = 70
$ws.Range("AM9").Value = 75
$ws.Range("AN9").Value = 3.2
$ws.Range("AO9").Value = 7.2
$ws.Range("AP9").Value = 19.5
$ws.Range("AQ9").Value = 23
$ws.Range("AR9").Value = 60
$ws.Range("AS9").Value = 300
$ws.Range("AT9").Value = 2.67
$ws.Range("AU9").Value = 8.75
$ws.Range("AV9").Value = 100
$ws.Range("AW9").Value = 51
$ws.Range("AX9").Value = 7.2
$ws.Range("AY9").Value = 37
$ws.Range("AZ9").Value = 45
$ws.Range("BA9").Value = 250
$ws.Range("BB9").Value = 300
$ws.Range("BC9").Value = 500
$ws.Range("BD9").Value = 51

# --- Small odds tweaks on the Ukraine match, now at row 10 ---
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 7.33

Write-Output "edit complete"
